# Mødeprotokol - attendance update for the week of 2017-04-18 .. 2017-04-24.
# The sheet pre-fills upcoming days with "Ikke registreret" (not yet
# registered); as each day passes, the real attendance status is typed in.
# This continues filling in the days 2017-04-19 .. 2017-04-24 (rows 79-84)
# and corrects a stray/garbled entry at I28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fremmøde Stats")

# Fix a previously mis-typed attendance note - should simply read "Mødt".
$ws.Range("I28").Value = "Mødt"

# Wednesday 2017-04-19 - everyone showed up, Mustafa called in sick.
$ws.Range("C79:H79").Value = "Mødt"
$ws.Range("I79").Value = "Kom ikke / Syg"

# Thursday 2017-04-20 - same as the day before.
$ws.Range("C80:H80").Value = "Mødt"
$ws.Range("I80").Value = "Kom ikke / Syg"

# Friday 2017-04-21 - everyone was out on agreed leave.
$ws.Range("C81:I81").Value = "Aftalt"

# Monday 2017-04-24 - Casper S called in sick, Mustafa was running late,
# everyone else showed up.
$ws.Range("C84").Value = "Mødt"
$ws.Range("D84").Value = "Kom ikke / Syg"
$ws.Range("E84:H84").Value = "Mødt"
$ws.Range("I84").Value = "Forsinket / 11:40"

# Keep the viewport/selection in sync with the last-edited cell.
$ws.Activate() | Out-Null
$ws.Range("I84").Select() | Out-Null
